$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 6667266
$ws.Range("I12").Value = 11111267
$ws.Range("K12").Value = 11111267
$ws.Range("M12").Value = -11111097
$ws.Range("H17").Value = 6668398
$ws.Range("J17").Value = 7144641
$ws.Range("L17").Value = 21433923
$ws.Range("N17").Value = -21434259
$ws.Range("H28").Value = 991.9167
$ws.Range("I28").Value = 912
$ws.Range("K28").Value = 912
$ws.Range("M28").Value = -427
$ws.Range("H62").Value = 1612.5
$ws.Range("I62").Value = 1271.4286
$ws.Range("K62").Value = 1271.4286
$ws.Range("M62").Value = -647.4286
$ws.Range("H64").Value = 9941.933999999999
$ws.Range("I64").Value = 3532.75
$ws.Range("J64").Value = 12272.546
$ws.Range("K64").Value = 3532.75
$ws.Range("L64").Value = 12272.546
$ws.Range("M64").Value = -3284.75
$ws.Range("N64").Value = -12768.546
$ws.Range("H65").Value = 1612.5
$ws.Range("I65").Value = 1271.4286
$ws.Range("K65").Value = 6357.143
$ws.Range("M65").Value = -3237.143
$ws.Range("H67").Value = 9941.933999999999
$ws.Range("I67").Value = 3532.75
$ws.Range("J67").Value = 12272.546
$ws.Range("K67").Value = 3532.75
$ws.Range("L67").Value = 12272.546
$ws.Range("M67").Value = -2674.75
$ws.Range("N67").Value = -13988.546
$ws.Range("H98").Value = 861.8276
$ws.Range("I98").Value = 869.1539
$ws.Range("K98").Value = 869.1539
$ws.Range("M98").Value = 628.8461
$ws.Range("H107").Value = 676.7727
$ws.Range("I107").Value = 755.8421
$ws.Range("K107").Value = 755.8421
$ws.Range("M107").Value = 1164.1579
$ws.Range("H122").Value = 861.8276
$ws.Range("I122").Value = 869.1539
$ws.Range("K122").Value = 2607.4617
$ws.Range("M122").Value = -157.4616999999998
$ws.Range("H138").Value = 2029
$ws.Range("I138").Value = 1378.1666
$ws.Range("J138").Value = 2321.875
$ws.Range("K138").Value = 4134.4998
$ws.Range("L138").Value = 6965.625
$ws.Range("M138").Value = 1005.5002
$ws.Range("N138").Value = -17245.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1902.963
$ws.Range("I2").Value = 908.55
$ws.Range("K2").Value = 908.55
$ws.Range("M2").Value = -795.55
$ws.Range("H32").Value = 2061.7827
$ws.Range("I32").Value = 2117
$ws.Range("K32").Value = 2117
$ws.Range("M32").Value = -1830
$ws.Range("H45").Value = 3556.6
$ws.Range("I45").Value = 1963.5454
$ws.Range("J45").Value = 7937.5
$ws.Range("K45").Value = 1963.5454
$ws.Range("L45").Value = 7937.5
$ws.Range("M45").Value = -1586.5454
$ws.Range("N45").Value = -8691.5
$ws.Range("H116").Value = 1902.963
$ws.Range("I116").Value = 908.55
$ws.Range("K116").Value = 908.55
$ws.Range("M116").Value = 1385.45
$ws.Range("H122").Value = 2764.16
$ws.Range("I122").Value = 2474.4375
$ws.Range("J122").Value = 3279.2222
$ws.Range("K122").Value = 7423.3125
$ws.Range("L122").Value = 9837.6666
$ws.Range("M122").Value = -4973.3125
$ws.Range("N122").Value = -14737.6666
$ws.Range("H132").Value = 4601
$ws.Range("I132").Value = 3759.818
$ws.Range("K132").Value = 11279.454
$ws.Range("M132").Value = -8749.454000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1902.963
$ws.Range("I3").Value = 908.55
$ws.Range("K3").Value = 908.55
$ws.Range("M3").Value = -794.55
$ws.Range("H86").Value = 1626.8182
$ws.Range("I86").Value = 1589.5
$ws.Range("K86").Value = 1589.5
$ws.Range("M86").Value = -466.5
$ws.Range("H89").Value = 1626.8182
$ws.Range("I89").Value = 1589.5
$ws.Range("K89").Value = 7947.5
$ws.Range("M89").Value = -2331.5
$ws.Range("H107").Value = 1757.0834
$ws.Range("I107").Value = 1718.0952
$ws.Range("K107").Value = 1718.0952
$ws.Range("M107").Value = 201.9048

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 80763.42999999999
$ws.Range("I86").Value = 10481.667
$ws.Range("J86").Value = 133474.75
$ws.Range("K86").Value = 10481.667
$ws.Range("L86").Value = 133474.75
$ws.Range("M86").Value = -9358.666999999999
$ws.Range("N86").Value = -135720.75
$ws.Range("H89").Value = 80763.42999999999
$ws.Range("I89").Value = 10481.667
$ws.Range("J89").Value = 133474.75
$ws.Range("K89").Value = 52408.335
$ws.Range("L89").Value = 667373.75
$ws.Range("M89").Value = -46792.335
$ws.Range("N89").Value = -678605.75
$ws.Range("H99").Value = 4428.25
$ws.Range("I99").Value = 3954.3333
$ws.Range("J99").Value = 5850
$ws.Range("K99").Value = 3954.3333
$ws.Range("L99").Value = 5850
$ws.Range("M99").Value = -2456.3333
$ws.Range("N99").Value = -8846
$ws.Range("H107").Value = 854.7778
$ws.Range("I107").Value = 824.75
$ws.Range("J107").Value = 878.8
$ws.Range("K107").Value = 824.75
$ws.Range("L107").Value = 878.8
$ws.Range("M107").Value = 1095.25
$ws.Range("N107").Value = -4718.8
$ws.Range("H122").Value = 3888.7
$ws.Range("I122").Value = 3251.4119
$ws.Range("K122").Value = 9754.235700000001
$ws.Range("M122").Value = -7304.235700000001
$ws.Range("H126").Value = 4428.25
$ws.Range("I126").Value = 3954.3333
$ws.Range("J126").Value = 5850
$ws.Range("K126").Value = 11862.9999
$ws.Range("L126").Value = 17550
$ws.Range("M126").Value = -9392.999899999999
$ws.Range("N126").Value = -22490
$ws.Range("H132").Value = 2401.111
$ws.Range("I132").Value = 2087.9583
$ws.Range("K132").Value = 6263.874899999999
$ws.Range("M132").Value = -3733.874899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 11297.111
$ws.Range("I23").Value = 210.71428
$ws.Range("K23").Value = 632.14284
$ws.Range("M23").Value = -397.14284
$ws.Range("H46").Value = 547.5
$ws.Range("I46").Value = 590
$ws.Range("J46").Value = 533.3333
$ws.Range("K46").Value = 1770
$ws.Range("L46").Value = 1599.9999
$ws.Range("M46").Value = -1679
$ws.Range("N46").Value = -1781.9999
$ws.Range("H121").Value = 59497.06
$ws.Range("I121").Value = 533.2222
$ws.Range("K121").Value = 1599.6666
$ws.Range("M121").Value = -289.6666
$ws.Range("H132").Value = 1716.9333
$ws.Range("I132").Value = 1545.4
$ws.Range("K132").Value = 13908.6
$ws.Range("M132").Value = -11378.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1374.5
$ws.Range("I31").Value = 1374.5
$ws.Range("K31").Value = 1374.5
$ws.Range("M31").Value = -1082.5
$ws.Range("H37").Value = 1374.5
$ws.Range("I37").Value = 1374.5
$ws.Range("K37").Value = 1374.5
$ws.Range("M37").Value = -1097.5
$ws.Range("H80").Value = 6336.8335
$ws.Range("I80").Value = 5717.875
$ws.Range("J80").Value = 7574.75
$ws.Range("K80").Value = 5717.875
$ws.Range("L80").Value = 7574.75
$ws.Range("M80").Value = -4719.875
$ws.Range("N80").Value = -9570.75
$ws.Range("H83").Value = 6336.8335
$ws.Range("I83").Value = 5717.875
$ws.Range("J83").Value = 7574.75
$ws.Range("K83").Value = 28589.375
$ws.Range("L83").Value = 37873.75
$ws.Range("M83").Value = -23597.375
$ws.Range("N83").Value = -47857.75
$ws.Range("H97").Value = 882.2222
$ws.Range("I97").Value = 730.3333
$ws.Range("J97").Value = 1186
$ws.Range("K97").Value = 730.3333
$ws.Range("L97").Value = 1186
$ws.Range("M97").Value = -234.3333
$ws.Range("N97").Value = -2178
$ws.Range("H102").Value = 2448.875
$ws.Range("I102").Value = 1763
$ws.Range("J102").Value = 7250
$ws.Range("K102").Value = 1763
$ws.Range("L102").Value = 7250
$ws.Range("M102").Value = -141
$ws.Range("N102").Value = -10494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3999.6667
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3999.6667
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3999.6667
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4223.6667
$ws.Range("H40").Value = 2454.4285
$ws.Range("I40").Value = 2446.8333
$ws.Range("K40").Value = 2446.8333
$ws.Range("M40").Value = -2310.8333
$ws.Range("H43").Value = 6000
$ws.Range("I43").Value = 6000
$ws.Range("K43").Value = 6000
$ws.Range("M43").Value = -5807
$ws.Range("H122").Value = 4247.7144
$ws.Range("I122").Value = 3911.6667
$ws.Range("K122").Value = 11735.0001
$ws.Range("M122").Value = -9285.000100000001
$ws.Range("H126").Value = 3999.6667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3999.6667
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11999.0001
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16939.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9497.5
$ws.Range("J15").Value = 9495
$ws.Range("L15").Value = 9495
$ws.Range("N15").Value = -10071
$ws.Range("H107").Value = 629
$ws.Range("I107").Value = 413.8
$ws.Range("J107").Value = 1167
$ws.Range("K107").Value = 1241.4
$ws.Range("L107").Value = 3501
$ws.Range("M107").Value = 678.5999999999999
$ws.Range("N107").Value = -7341
$ws.Range("H126").Value = 6955.1
$ws.Range("I126").Value = 6505.75
$ws.Range("K126").Value = 19517.25
$ws.Range("M126").Value = -17047.25
$ws.Range("H132").Value = 8333.333000000001
$ws.Range("I132").Value = 7500
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -19970
$ws.Range("N132").Value = -35060
